$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the image-dimension labels in column A to more "favourable" (power-of-two / square) sizes.
$ws.Range("A4").Value = "512x512"
$ws.Range("A6").Value = "2048x2048"
$ws.Range("A7").Value = "4096x4096"

# Reselect A9 (matches the selection recorded in the edited workbook)
$ws.Range("A9").Select()
